$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.658.76"
$ws.Range("E2").Value = "  -2.02%  "
$ws.Range("D3").Value = "3.389.99"
$ws.Range("E3").Value = "  -2.11%  "
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "405.86"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -2.03%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "133.33"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +8.56%  "
$ws.Range("E7").Value = "  -0.60%  "
$ws.Range("E9").Value = "  -2.12%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.121"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -6.09%  "
$ws.Range("E11").Value = "  +2.99%  "
$ws.Range("E12").Value = "  -1.32%  "
$ws.Range("D13").Value = "3.922.76"
$ws.Range("E13").Value = "  -2.41%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "8.42"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -2.08%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "19.77"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -0.88%  "
$ws.Range("D16").Value = "3.408.71"
$ws.Range("E16").Value = "  -1.18%  "
$ws.Range("D17").Value = "61.685.29"
$ws.Range("E17").Value = "  -1.95%  "
$ws.Range("E18").Value = "  -1.28%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.99"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +1.00%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0000128"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -6.33%  "
$ws.Range("E21").Value = "  -3.65%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "85.33"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +4.62%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "317.85"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +0.55%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.77"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -0.58%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.12"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -1.36%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "4.78"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +10.87%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "29.55"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -4.10%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.30"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +6.08%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.68"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -0.99%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.69"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +3.48%  "
$ws.Range("E31").Value = "  -1.65%  "
$ws.Range("E32").Value = "  -0.34%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "11.39"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -1.80%  "
$ws.Range("E34").Value = "  -0.61%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "41.67"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -0.90%  "
$ws.Range("E36").Value = "  -2.17%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "51.76"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -0.80%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.999"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +0.09%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.43"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -1.85%  "
$ws.Range("E40").Value = "  -2.94%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "139.34"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +2.70%  "
$ws.Range("E42").Value = "  -0.46%  "
$ws.Range("E43").Value = "  -1.03%  "
$ws.Range("E44").Value = "  +4.28%  "
$ws.Range("E45").Value = "  +2.67%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "16.69"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -0.79%  "
$ws.Range("E47").Value = "  -1.81%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "21.40"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -2.04%  "
$ws.Range("D49").Value = "2.126.71"
$ws.Range("E49").Value = "  -2.50%  "
$ws.Range("E50").Value = "  -7.00%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.89"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +1.60%  "
